# Apply updated cryptocurrency price (column D) and hourly volume change (column E)
# values, matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.850.86'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.09%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.223.81'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.21%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '241.67'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.35%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.83%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '72.77'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.66%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.585'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.03%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.52'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.81%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0945'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.86%  '

$ws.Range('E12').Value = '  +0.61%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.88'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.64%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.551.42'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.38%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.21'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.77%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.830'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.58%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.211.81'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.66%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '41.710.87'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.99%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000105'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -4.74%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.17'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.03%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.90'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.33%  '

$ws.Range('E22').Value = '  +9.81%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '228.76'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.14%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.03'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -6.10%  '

$ws.Range('E25').Value = '  +0.04%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.32'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -4.01%  '

$ws.Range('E27').Value = '  -1.42%  '

$ws.Range('E28').Value = '  -2.08%  '

$ws.Range('E29').Value = '  -1.09%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '166.89'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.08%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.43'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.11%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0793'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.68%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.49'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.43%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '30.76'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.14%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.124'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.29%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.109'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -7.50%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.26'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.80%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0302'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.61%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '13.07'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.36%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.11'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.60%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '64.10'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.21%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.61'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.42%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.196'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.64%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.65'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.28%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '102.21'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.65%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0994'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.42%  '

$ws.Range('E47').Value = '  -0.78%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.17'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.09%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.32'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -3.34%  '

$ws.Range('E50').Value = '  -1.52%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.425.17'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.52%  '
